# Auto-generated script applying updated Leve profit values
# per the commit's scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 376.125
$ws.Cells.Item(28, 9).Value = 339.41666
$ws.Cells.Item(28, 10).Value = 486.25
$ws.Cells.Item(28, 11).Value = 339.41666
$ws.Cells.Item(28, 12).Value = 486.25
$ws.Cells.Item(28, 13).Value = 145.58334
$ws.Cells.Item(28, 14).Value = -1456.25
$ws.Cells.Item(32, 8).Value = 166667260
$ws.Cells.Item(32, 10).Value = 726.6667
$ws.Cells.Item(32, 12).Value = 726.6667
$ws.Cells.Item(32, 14).Value = -1378.6667
$ws.Cells.Item(98, 8).Value = 650
$ws.Cells.Item(98, 9).Value = 600.36365
$ws.Cells.Item(98, 10).Value = 832
$ws.Cells.Item(98, 11).Value = 600.36365
$ws.Cells.Item(98, 12).Value = 832
$ws.Cells.Item(98, 13).Value = 897.63635
$ws.Cells.Item(98, 14).Value = -3828
$ws.Cells.Item(99, 8).Value = 2219.875
$ws.Cells.Item(99, 9).Value = 244.5
$ws.Cells.Item(99, 10).Value = 4195.25
$ws.Cells.Item(99, 11).Value = 733.5
$ws.Cells.Item(99, 12).Value = 12585.75
$ws.Cells.Item(99, 13).Value = 764.5
$ws.Cells.Item(99, 14).Value = -15581.75
$ws.Cells.Item(107, 8).Value = 200322.2
$ws.Cells.Item(107, 9).Value = 250151.25
$ws.Cells.Item(107, 10).Value = 1006
$ws.Cells.Item(107, 11).Value = 250151.25
$ws.Cells.Item(107, 12).Value = 1006
$ws.Cells.Item(107, 13).Value = -248231.25
$ws.Cells.Item(107, 14).Value = -4846
$ws.Cells.Item(111, 8).Value = 6009.4443
$ws.Cells.Item(111, 9).Value = 1110.5
$ws.Cells.Item(111, 10).Value = 7409.143
$ws.Cells.Item(111, 11).Value = 3331.5
$ws.Cells.Item(111, 12).Value = 22227.429
$ws.Cells.Item(111, 13).Value = -264.5
$ws.Cells.Item(111, 14).Value = -28361.429
$ws.Cells.Item(122, 8).Value = 650
$ws.Cells.Item(122, 9).Value = 600.36365
$ws.Cells.Item(122, 10).Value = 832
$ws.Cells.Item(122, 11).Value = 1801.09095
$ws.Cells.Item(122, 12).Value = 2496
$ws.Cells.Item(122, 13).Value = 648.90905
$ws.Cells.Item(122, 14).Value = -7396
$ws.Cells.Item(125, 8).Value = 1868.1578
$ws.Cells.Item(125, 9).Value = 668.4286
$ws.Cells.Item(125, 10).Value = 2568
$ws.Cells.Item(125, 11).Value = 6015.8574
$ws.Cells.Item(125, 12).Value = 23112
$ws.Cells.Item(125, 13).Value = -3555.8574
$ws.Cells.Item(125, 14).Value = -28032
$ws.Cells.Item(132, 8).Value = 4826.6724
$ws.Cells.Item(132, 9).Value = 2967.8948
$ws.Cells.Item(132, 10).Value = 8358.35
$ws.Cells.Item(132, 11).Value = 8903.6844
$ws.Cells.Item(132, 12).Value = 25075.05
$ws.Cells.Item(132, 13).Value = -6373.6844
$ws.Cells.Item(132, 14).Value = -30135.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 3000
$ws.Cells.Item(97, 9).Value = 1933.3334
$ws.Cells.Item(97, 10).Value = 6200
$ws.Cells.Item(97, 11).Value = 1933.3334
$ws.Cells.Item(97, 12).Value = 6200
$ws.Cells.Item(97, 13).Value = -1437.3334
$ws.Cells.Item(97, 14).Value = -7192
$ws.Cells.Item(102, 8).Value = 52633156
$ws.Cells.Item(102, 9).Value = 1601.2142
$ws.Cells.Item(102, 10).Value = 200001500
$ws.Cells.Item(102, 11).Value = 1601.2142
$ws.Cells.Item(102, 12).Value = 200001500
$ws.Cells.Item(102, 13).Value = 20.78580000000011
$ws.Cells.Item(102, 14).Value = -200004744

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(34, 8).Value = 16975.334
$ws.Cells.Item(34, 10).Value = 16975.334
$ws.Cells.Item(34, 12).Value = 16975.334
$ws.Cells.Item(34, 14).Value = -17203.334
$ws.Cells.Item(86, 8).Value = 2294.3713
$ws.Cells.Item(86, 9).Value = 2406.7585
$ws.Cells.Item(86, 11).Value = 2406.7585
$ws.Cells.Item(86, 13).Value = -1283.7585
$ws.Cells.Item(89, 8).Value = 2294.3713
$ws.Cells.Item(89, 9).Value = 2406.7585
$ws.Cells.Item(89, 11).Value = 12033.7925
$ws.Cells.Item(89, 13).Value = -6417.7925
$ws.Cells.Item(105, 8).Value = 3065.9412
$ws.Cells.Item(105, 9).Value = 1901.25
$ws.Cells.Item(105, 10).Value = 4101.222
$ws.Cells.Item(105, 11).Value = 1901.25
$ws.Cells.Item(105, 12).Value = 4101.222
$ws.Cells.Item(105, 13).Value = -154.25
$ws.Cells.Item(105, 14).Value = -7595.222
$ws.Cells.Item(134, 8).Value = 1418.2821
$ws.Cells.Item(134, 9).Value = 1200.7742
$ws.Cells.Item(134, 10).Value = 2261.125
$ws.Cells.Item(134, 11).Value = 3602.3226
$ws.Cells.Item(134, 12).Value = 6783.375
$ws.Cells.Item(134, 13).Value = -1067.3226
$ws.Cells.Item(134, 14).Value = -11853.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 965.2917
$ws.Cells.Item(58, 9).Value = 1042.7894
$ws.Cells.Item(58, 10).Value = 670.8
$ws.Cells.Item(58, 11).Value = 1042.7894
$ws.Cells.Item(58, 12).Value = 670.8
$ws.Cells.Item(58, 13).Value = -839.7893999999999
$ws.Cells.Item(58, 14).Value = -1076.8
$ws.Cells.Item(62, 8).Value = 2733.8823
$ws.Cells.Item(62, 9).Value = 2738.4666
$ws.Cells.Item(62, 10).Value = 2699.5
$ws.Cells.Item(62, 11).Value = 2738.4666
$ws.Cells.Item(62, 12).Value = 2699.5
$ws.Cells.Item(62, 13).Value = -2114.4666
$ws.Cells.Item(62, 14).Value = -3947.5
$ws.Cells.Item(65, 8).Value = 2733.8823
$ws.Cells.Item(65, 9).Value = 2738.4666
$ws.Cells.Item(65, 10).Value = 2699.5
$ws.Cells.Item(65, 11).Value = 13692.333
$ws.Cells.Item(65, 12).Value = 13497.5
$ws.Cells.Item(65, 13).Value = -10572.333
$ws.Cells.Item(65, 14).Value = -19737.5
$ws.Cells.Item(86, 8).Value = 7429.2
$ws.Cells.Item(86, 9).Value = 11344.75
$ws.Cells.Item(86, 10).Value = 2954.2856
$ws.Cells.Item(86, 11).Value = 11344.75
$ws.Cells.Item(86, 12).Value = 2954.2856
$ws.Cells.Item(86, 13).Value = -10221.75
$ws.Cells.Item(86, 14).Value = -5200.2856
$ws.Cells.Item(89, 8).Value = 7429.2
$ws.Cells.Item(89, 9).Value = 11344.75
$ws.Cells.Item(89, 10).Value = 2954.2856
$ws.Cells.Item(89, 11).Value = 56723.75
$ws.Cells.Item(89, 12).Value = 14771.428
$ws.Cells.Item(89, 13).Value = -51107.75
$ws.Cells.Item(89, 14).Value = -26003.428
$ws.Cells.Item(111, 8).Value = 37801.332
$ws.Cells.Item(111, 10).Value = 37801.332
$ws.Cells.Item(111, 12).Value = 37801.332
$ws.Cells.Item(111, 14).Value = -45981.332
$ws.Cells.Item(136, 8).Value = 965.2917
$ws.Cells.Item(136, 9).Value = 1042.7894
$ws.Cells.Item(136, 10).Value = 670.8
$ws.Cells.Item(136, 11).Value = 3128.3682
$ws.Cells.Item(136, 12).Value = 2012.4
$ws.Cells.Item(136, 13).Value = -578.3681999999999
$ws.Cells.Item(136, 14).Value = -7112.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 43515
$ws.Cells.Item(2, 9).Value = 142883.14
$ws.Cells.Item(2, 10).Value = 41.4375
$ws.Cells.Item(2, 11).Value = 857298.8400000001
$ws.Cells.Item(2, 12).Value = 248.625
$ws.Cells.Item(2, 13).Value = -857185.8400000001
$ws.Cells.Item(2, 14).Value = -474.625
$ws.Cells.Item(40, 8).Value = 2096613.9
$ws.Cells.Item(40, 9).Value = 171.57143
$ws.Cells.Item(40, 11).Value = 686.28572
$ws.Cells.Item(40, 13).Value = -617.28572
$ws.Cells.Item(92, 8).Value = 190
$ws.Cells.Item(92, 9).Value = 180
$ws.Cells.Item(92, 11).Value = 540
$ws.Cells.Item(92, 13).Value = 708
$ws.Cells.Item(107, 8).Value = 604.29114
$ws.Cells.Item(107, 9).Value = 333.23914
$ws.Cells.Item(107, 10).Value = 982.1212
$ws.Cells.Item(107, 11).Value = 999.7174200000001
$ws.Cells.Item(107, 12).Value = 2946.3636
$ws.Cells.Item(107, 13).Value = 920.2825799999999
$ws.Cells.Item(107, 14).Value = -6786.363600000001
$ws.Cells.Item(113, 8).Value = 530
$ws.Cells.Item(113, 9).Value = 700.5
$ws.Cells.Item(113, 10).Value = 487.375
$ws.Cells.Item(113, 11).Value = 2101.5
$ws.Cells.Item(113, 12).Value = 1462.125
$ws.Cells.Item(113, 13).Value = 68.5
$ws.Cells.Item(113, 14).Value = -5802.125
$ws.Cells.Item(117, 8).Value = 2448.1
$ws.Cells.Item(117, 9).Value = 296.2
$ws.Cells.Item(117, 10).Value = 4600
$ws.Cells.Item(117, 11).Value = 888.5999999999999
$ws.Cells.Item(117, 12).Value = 13800
$ws.Cells.Item(117, 13).Value = 2553.4
$ws.Cells.Item(117, 14).Value = -20684

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 7340.6665
$ws.Cells.Item(97, 9).Value = 2500
$ws.Cells.Item(97, 10).Value = 9761
$ws.Cells.Item(97, 11).Value = 2500
$ws.Cells.Item(97, 12).Value = 9761
$ws.Cells.Item(97, 13).Value = -2004
$ws.Cells.Item(97, 14).Value = -10753
$ws.Cells.Item(108, 8).Value = 38400
$ws.Cells.Item(108, 10).Value = 38400
$ws.Cells.Item(108, 12).Value = 38400
$ws.Cells.Item(108, 14).Value = -46080

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 814.3570999999999
$ws.Cells.Item(22, 9).Value = 657
$ws.Cells.Item(22, 11).Value = 657
$ws.Cells.Item(22, 13).Value = -362
$ws.Cells.Item(27, 8).Value = 814.3570999999999
$ws.Cells.Item(27, 9).Value = 657
$ws.Cells.Item(27, 11).Value = 657
$ws.Cells.Item(27, 13).Value = -550
$ws.Cells.Item(132, 8).Value = 2469.8572
$ws.Cells.Item(132, 9).Value = 1870.4546
$ws.Cells.Item(132, 10).Value = 4667.6665
$ws.Cells.Item(132, 11).Value = 5611.3638
$ws.Cells.Item(132, 12).Value = 14002.9995
$ws.Cells.Item(132, 13).Value = -3081.3638
$ws.Cells.Item(132, 14).Value = -19062.9995
$ws.Cells.Item(136, 8).Value = 1652
$ws.Cells.Item(136, 9).Value = 1526.6842
$ws.Cells.Item(136, 10).Value = 2842.5
$ws.Cells.Item(136, 11).Value = 4580.0526
$ws.Cells.Item(136, 12).Value = 8527.5
$ws.Cells.Item(136, 13).Value = -2030.0526
$ws.Cells.Item(136, 14).Value = -13627.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 253409.25
$ws.Cells.Item(113, 9).Value = 500418.5
$ws.Cells.Item(113, 10).Value = 6400
$ws.Cells.Item(113, 11).Value = 1501255.5
$ws.Cells.Item(113, 12).Value = 19200
$ws.Cells.Item(113, 13).Value = -1499085.5
$ws.Cells.Item(113, 14).Value = -23540
